{"js": "// Replace the multiplication-problem text in each table cell per the\n// old-value -> new-value mapping derived from the OOXML diff. All 25\n// old values are unique within the document, so a single search+replace\n// pass per pair is unambiguous.\nconst replacements = [\n  [\"14\u00d757=\", \"39\u00d793=\"],\n  [\"73\u00d757=\", \"83\u00d771=\"],\n  [\"31\u00d762=\", \"77\u00d784=\"],\n  [\"80\u00d720=\", \"38\u00d737=\"],\n  [\"77\u00d760=\", \"34\u00d716=\"],\n  [\"85\u00d796=\", \"98\u00d764=\"],\n  [\"16\u00d724=\", \"38\u00d752=\"],\n  [\"92\u00d717=\", \"51\u00d770=\"],\n  [\"76\u00d714=\", \"60\u00d754=\"],\n  [\"87\u00d789=\", \"46\u00d775=\"],\n  [\"14\u00d762=\", \"88\u00d781=\"],\n  [\"43\u00d799=\", \"72\u00d779=\"],\n  [\"65\u00d767=\", \"85\u00d764=\"],\n  [\"93\u00d788=\", \"57\u00d746=\"],\n  [\"65\u00d757=\", \"55\u00d777=\"],\n  [\"93\u00d730=\", \"74\u00d754=\"],\n  [\"40\u00d714=\", \"49\u00d790=\"],\n  [\"65\u00d747=\", \"60\u00d769=\"],\n  [\"70\u00d739=\", \"85\u00d711=\"],\n  [\"90\u00d763=\", \"59\u00d783=\"],\n  [\"14\u00d714=\", \"43\u00d771=\"],\n  [\"57\u00d787=\", \"97\u00d795=\"],\n  [\"63\u00d749=\", \"75\u00d742=\"],\n  [\"35\u00d737=\", \"60\u00d746=\"],\n  [\"73\u00d770=\", \"60\u00d712=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the multiplication-problem text in each table cell per the\n# old-value -> new-value mapping derived from the OOXML diff. All 25\n# old values are unique within the document, so Find/Replace (ReplaceAll)\n# for each pair is unambiguous and touches exactly one cell.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"14\u00d757=\", \"39\u00d793=\"),\n    @(\"73\u00d757=\", \"83\u00d771=\"),\n    @(\"31\u00d762=\", \"77\u00d784=\"),\n    @(\"80\u00d720=\", \"38\u00d737=\"),\n    @(\"77\u00d760=\", \"34\u00d716=\"),\n    @(\"85\u00d796=\", \"98\u00d764=\"),\n    @(\"16\u00d724=\", \"38\u00d752=\"),\n    @(\"92\u00d717=\", \"51\u00d770=\"),\n    @(\"76\u00d714=\", \"60\u00d754=\"),\n    @(\"87\u00d789=\", \"46\u00d775=\"),\n    @(\"14\u00d762=\", \"88\u00d781=\"),\n    @(\"43\u00d799=\", \"72\u00d779=\"),\n    @(\"65\u00d767=\", \"85\u00d764=\"),\n    @(\"93\u00d788=\", \"57\u00d746=\"),\n    @(\"65\u00d757=\", \"55\u00d777=\"),\n    @(\"93\u00d730=\", \"74\u00d754=\"),\n    @(\"40\u00d714=\", \"49\u00d790=\"),\n    @(\"65\u00d747=\", \"60\u00d769=\"),\n    @(\"70\u00d739=\", \"85\u00d711=\"),\n    @(\"90\u00d763=\", \"59\u00d783=\"),\n    @(\"14\u00d714=\", \"43\u00d771=\"),\n    @(\"57\u00d787=\", \"97\u00d795=\"),\n    @(\"63\u00d749=\", \"75\u00d742=\"),\n    @(\"35\u00d737=\", \"60\u00d746=\"),\n    @(\"73\u00d770=\", \"60\u00d712=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
